$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new training-day column (AS) is appended after the existing last day
# column (AR). Copy AR's per-row formatting into AS first (so the new
# column picks up the same style indices as the rest of the table), then
# fill in the actual values.

$ws.Range("AR1:AR29").Copy()
$ws.Range("AS1:AS29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row: new session date (10 Sept 2025 -> serial 45910)
$ws.Range("AS1").Value2 = 45910

# Attendance marks for the new session, one per player row.
$attendance = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "B"
    11 = "P"
    12 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "A"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $attendance.Keys) {
    $ws.Cells.Item($row, 45).Value2 = $attendance[$row]
}

# Keep the active selection one column further right, same as the
# author's view after filling in the new column.
$ws.Range("AU27").Select()
